# Add a "Save" column (column H) to the s_vals sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style from the existing "sum" header cell (G1) onto the
# new H1 header cell, then overwrite the copied value with "Save".
$ws.Cells.Item(1, 7).Copy($ws.Cells.Item(1, 8))
$ws.Cells.Item(1, 8).Value = "Save"

# Fill in the "Save" values for each data row.
$ws.Cells.Item(2, 8).Value = 0
$ws.Cells.Item(3, 8).Value = 1
